$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (copy H1's formatting/style so I1/J1 match the bold,
# centered, bordered header style already used by the other headers)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-13
$values = @{
    2  = @(9, 9)
    3  = @(10, 10)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(8, 8)
    7  = @(7, 7)
    8  = @(8, 8)
    9  = @(5, 7)
    10 = @(6, 6)
    11 = @(6, 8)
    12 = @(6, 6)
    13 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
